$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Encoded data: patientId|percentage pairs, comma separated, in target row order (A2:B382)
$data = "369|93.38,278|85.14,162|84.27,325|83.84,208|82.86,160|82.75,166|82.65,268|82.65,267|82,161|80.91,377|80.8,363|80.69,335|80.37,328|79.93,287|79.61,326|79.61,298|79.5,301|79.5,320|79.28,303|79.18,209|79.07,157|78.96,366|78.74,312|78.52,274|78.42,273|78.2,302|78.09,323|78.09,271|77.98,358|77.98,283|77.87,154|77.77,333|77.66,352|77.66,262|77.55,276|77.55,289|77.55,308|77.55,259|77.22,294|77.22,304|77.11,334|77.11,306|76.9,329|76.79,285|76.68,293|76.68,153|76.57,338|76.46,349|75.7,381|75.7,292|75.6,300|75.6,305|75.6,337|75.49,357|75.38,272|75.27,288|75.27,332|75.27,346|75.16,286|75.05,296|75.05,345|75.05,200|74.95,313|74.95,317|74.95,281|74.84,307|74.84,277|74.73,282|74.73,275|74.62,344|74.62,260|74.4,330|74.4,197|74.3,379|74.3,336|74.19,311|73.97,347|73.86,365|73.86,373|73.75,170|73.64,324|73.54,340|73.54,375|73.54,165|73.21,370|73.21,211|73.1,378|73.1,343|72.99,364|72.99,269|72.89,169|72.78,181|72.78,196|72.78,279|72.78,280|72.78,321|72.78,376|72.78,353|72.67,205|72.56,299|72.56,185|72.45,193|72.45,266|72.34,192|72.23,258|72.23,339|72.23,342|72.13,362|72.13,284|72.02,331|72.02,341|71.91,163|71.69,318|71.58,190|71.48,310|71.37,355|71.37,315|71.26,360|71.26,270|71.15,159|71.04,182|71.04,309|71.04,327|71.04,351|70.93,177|70.82,322|70.72,261|70.5,173|70.28,210|70.28,290|70.28,263|70.07,319|70.07,202|69.85,348|69.85,79|69.74,81|69.74,172|69.63,295|69.41,68|69.2,106|69.09,168|68.98,184|68.98,255|68.98,264|68.98,195|68.66,350|68.66,47|68.44,176|68.33,291|68.33,155|68.22,265|68.22,31|68.11,120|68.11,367|68.11,314|68,354|68,356|68,359|68,297|67.9,316|67.9,71|67.79,117|67.79,206|67.79,52|67.46,198|67.35,118|67.25,167|67.25,232|67.25,188|67.14,158|67.03,95|66.92,252|66.81,380|66.81,187|66.7,63|66.49,107|66.49,112|66.49,180|66.49,139|66.38,179|66.38,156|66.27,371|66.16,189|65.84,83|65.73,144|65.73,178|65.73,256|65.73,374|65.62,67|65.51,78|65.51,90|65.51,183|65.51,250|65.51,171|65.29,203|65.29,236|65.29,93|65.18,253|65.18,59|65.08,372|65.08,61|64.97,234|64.97,111|64.86,64|64.64,194|64.53,246|64.43,80|64.32,105|64.32,174|64.32,240|63.99,72|63.88,60|63.77,53|63.56,96|63.56,73|63.34,164|63.34,368|63.34,361|63.23,74|63.02,175|62.91,48|62.8,75|62.8,249|62.8,77|62.58,103|62.58,191|62.58,237|62.58,148|62.47,69|62.26,135|62.04,199|62.04,15|61.93,92|61.93,57|61.71,91|61.71,65|61.61,142|61.39,245|61.39,257|61.39,121|61.28,124|61.28,186|61.28,2|61.17,55|61.17,62|61.17,86|61.17,102|61.17,87|60.95,146|60.95,151|60.95,233|60.95,207|60.85,58|60.74,98|60.74,66|60.63,247|60.63,43|60.52,201|60.52,127|60.2,239|60.2,50|59.98,214|59.87,251|59.87,54|59.65,56|59.65,82|59.65,94|59.65,128|59.65,238|59.65,230|59.54,126|59.44,244|59.33,145|59.22,108|59.11,147|59.11,254|59,84|58.89,137|58.89,222|58.89,204|58.79,224|58.68,109|58.57,143|58.57,125|58.46,227|58.35,113|58.24,248|58.24,218|58.13,243|58.13,129|58.03,130|58.03,115|57.92,76|57.81,219|57.81,51|57.7,134|57.7,116|57.59,215|57.59,141|57.48,242|57.48,221|57.38,241|57.38,216|57.27,150|57.16,217|57.16,220|57.05,138|56.94,213|56.94,85|56.83,123|56.83,152|56.83,9|56.72,28|56.62,35|56.62,104|56.51,114|56.51,100|56.29,149|56.18,231|56.07,11|55.97,49|55.75,70|55.75,122|55.42,99|55.31,101|54.66,18|54.34,229|54.12,44|53.8,12|53.69,23|53.69,34|53.69,212|53.47,225|53.47,29|53.15,5|53.04,133|53.04,46|52.93,136|52.82,228|52.82,33|52.6,132|52.49,13|52.39,226|52.39,16|52.28,19|52.28,140|52.28,36|51.95,30|51.84,1|51.74,14|51.63,20|51.52,38|51.52,17|51.3,3|51.19,37|51.19,235|51.19,131|50.87,25|50.76,39|50.54,97|50.43,26|50.33,45|50.22,8|49.78,223|49.78,119|49.46,10|49.24,32|48.92,41|48.92,89|48.7,27|48.05,40|47.83,7|47.72,24|47.4,21|46.96,110|46.75,22|46.2,88|46.2,42|46.1,6|45.88,4|45.66"

$pairs = $data.Split(",")
$rowCount = $pairs.Length

$startRow = 2
for ($i = 0; $i -lt $rowCount; $i++) {
    $parts = $pairs[$i].Split("|")
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = [double]$parts[0]
    $ws.Cells.Item($r, 2).Value = [double]$parts[1]
}
